# Updated for new display-read wescheme lib
#
# The shape holding the four "(require wescheme/...)" lines gets touched:
#   - para 1: the lone tab-run after "(require wescheme/oJ1vcDo5qd)" grows
#     from "<TAB>" to "<TAB>  " (tab + two spaces)
#   - para 2: the library id changes (1Q1f9pSrg8 -> RWJy5EoNzk); retyping it
#     by hand makes PowerPoint record it as five separate runs:
#     "(", "require", " ", "wescheme", "/RWJy5EoNzk) "
#   - para 3 & 4: same "<TAB>" -> "<TAB>  " widening as para 1, but here the
#     tab is part of the same run as the "(require ...)" text, so the whole
#     run text is rewritten in place (no new run boundary appears)

$p = $ppt.ActivePresentation

# Find the shape across the deck instead of hard-coding slide/shape indices.
$shape = $null
for ($si = 1; $si -le $p.Slides.Count -and $shape -eq $null; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $cand = $sl.Shapes.Item($shi)
        if ($cand.HasTextFrame) {
            if ($cand.TextFrame.HasText) {
                if ($cand.TextFrame.TextRange.Text -clike "*(require wescheme*") {
                    $shape = $cand
                    break
                }
            }
        }
    }
}

$tr = $shape.TextFrame.TextRange

# --- Paragraph 1: "(require wescheme/oJ1vcDo5qd)<TAB>; konnagrafiikkakirjasto"
$para1 = $tr.Paragraphs(1, 1)
$tabRun = $tr.Characters($para1.Start + 29, 1)   # the lone "<TAB>" run
$tabRun.Text = "`t  "

# --- Paragraph 2: "(require wescheme/1Q1f9pSrg8)<TAB>; kayttoliittymakirjasto"
$para2 = $tr.Paragraphs(2, 1)
$base2 = $para2.Start
$tr.Characters($base2, 1).Text       = "("
$tr.Characters($base2 + 1, 7).Text   = "require"
$tr.Characters($base2 + 8, 1).Text   = " "
$tr.Characters($base2 + 9, 8).Text   = "wescheme"
$tr.Characters($base2 + 17, 13).Text = "/RWJy5EoNzk) "

# --- Paragraph 3: "(require wescheme/f08DD6x94M)<TAB>; matikan apukirjasto"
$para3 = $tr.Paragraphs(3, 1)
$tr.Characters($para3.Start, 30).Text = "(require wescheme/f08DD6x94M)`t  "

# --- Paragraph 4: "(require wescheme/2W8inC9p82)<TAB>; kuvaajat ja diagrammit "
$para4 = $tr.Paragraphs(4, 1)
$tr.Characters($para4.Start, 30).Text = "(require wescheme/2W8inC9p82)`t  "
